# boxtet_setup.pptx: update box/tetrahedron figure (slide 1)
#
# 1) Reposition a handful of existing shapes/groups that make up the
#    box-and-tetrahedron diagram.
# 2) Add six new labeled text boxes describing the P->v / P->e / e->v
#    incidence tables for the box (left pair of cubes) and the
#    tetrahedron (right shape).
#
# NOTE on numeric literals: this COM host marshals Shape.Left/Top/Width/
# Height (and AddTextbox's arguments) through a single-precision (f32)
# point value before converting to EMU (1 pt = 12700 EMU) and truncating,
# so a naive `targetEmu/12700` literal can land one EMU short after the
# f32 round-trip. The literals below were solved so that, after the f32
# round-trip, they truncate to the exact target EMU.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

function Get-ShapeById($slide, $id) {
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $sh = $slide.Shapes.Item($i)
        if ($sh.Id -eq $id) {
            return $sh
        }
    }
    return $null
}

# ---------------------------------------------------------------------
# 1) Move existing shapes / groups
# ---------------------------------------------------------------------

(Get-ShapeById $s 19).Left = 32.70909501811023   # Group 18 (box group 2) -> 415405 EMU
(Get-ShapeById $s 32).Top = 12.0000391   # Group 31 (tetrahedron group) -> 152400 EMU
(Get-ShapeById $s 45).Left = 197.2500382   # v1 -> 2505075 EMU
(Get-ShapeById $s 46).Left = 0.0   # v0 -> 0 EMU
(Get-ShapeById $s 47).Left = 278.2613830826772   # v3 -> 3533919 EMU
(Get-ShapeById $s 48).Left = 79.46932983858268   # v2 -> 1009260 EMU
(Get-ShapeById $s 49).Left = 197.2500382   # v5 -> 2505075 EMU
(Get-ShapeById $s 50).Left = 0.0   # v4 -> 0 EMU
(Get-ShapeById $s 51).Left = 278.2613830826772   # v7 -> 3533919 EMU
(Get-ShapeById $s 52).Left = 79.46932983858268   # v6 -> 1009260 EMU
(Get-ShapeById $s 65).Top = 114.85586547165354   # p0 -> 1458669 EMU
(Get-ShapeById $s 66).Top = 123.46334458661418   # p2 -> 1567984 EMU
(Get-ShapeById $s 67).Top = 114.85586547165354   # p1 -> 1458669 EMU
(Get-ShapeById $s 68).Top = 28.609803179527557   # p5 -> 363344 EMU
(Get-ShapeById $s 69).Top = 190.37814333622046   # p4 -> 2417802 EMU
(Get-ShapeById $s 70).Top = 85.25625987244095   # p3 -> 1082754 EMU

# ---------------------------------------------------------------------
# 2) Add new incidence-table text boxes
# ---------------------------------------------------------------------

function Add-IncidenceTextBox($slide, $left, $top, $width, $height, $wrap, $name, $lines) {
    $tb = $slide.Shapes.AddTextbox(1, $left, $top, $width, $height)
    $tb.Name = $name
    if ($wrap -eq "none") {
        $tb.TextFrame.WordWrap = 0
    }
    $tb.TextFrame.AutoSize = 1
    $tb.Fill.Visible = 0
    $tb.TextFrame.TextRange.Text = [string]::Join("`r", $lines)
    [void]$tb
}

# box P->v table
Add-IncidenceTextBox $s 322.9839019677165 338.8781433362205 83.33098218188977 159.9468841937008 "none" "TextBox 94" @(
    "P->v",
    "0: 0,2,4,6",
    "1: 1,3,5,7",
    "2: 0,1,4,5",
    "3: 2,3,6,7",
    "4: 0,1,2,3",
    "5: 4,5,6,7"
)

# box P->e table
Add-IncidenceTextBox $s 417.42839055669293 339.13217168425194 101.75917435826771 159.9468841937008 "none" "TextBox 95" @(
    "P->e",
    "0: 0,2,8,10",
    "1: 1,3,9,11",
    "2: 4,6,8,9",
    "3: 5,7,10,11",
    "4: 0,1,4,5",
    "5: 2,3,6,7"
)

# box e->v table
Add-IncidenceTextBox $s 536.6646728692913 248.66720585433072 65.02893443779529 290.81256108503936 "none" "TextBox 96" @(
    "e->v",
    "0: 0,2 ",
    "1: 1,3",
    "2: 4,6",
    "3: 5,7",
    "4: 0,1",
    "5: 2,3",
    "6: 4,5",
    "7: 6,7",
    "8: 0,4",
    "9: 1,5",
    "10: 2,6",
    "11: 3,7"
)

# tetrahedron P->v table
Add-IncidenceTextBox $s 379.21357730708667 296.8405151409449 69.5728721456693 116.32507707007873 "none" "TextBox 3" @(
    "P->v",
    "0: 1,0,3",
    "1: 2,1,3",
    "2: 0,2,3",
    "3: 0,2,1"
)

# tetrahedron P->e table
Add-IncidenceTextBox $s 457.3654785708661 303.8124847448819 69.5728721456693 116.32507707007873 "none" "TextBox 38" @(
    "P->e",
    "0: 0,3,4",
    "1: 1,4,5",
    "2: 2,3,5",
    "3: 0,1,2"
)

# tetrahedron e->v table
Add-IncidenceTextBox $s 543.3597412393701 296.7219238637795 69.57578663149606 159.9468841937008 "square" "TextBox 39" @(
    "e->v",
    "0: 0,1",
    "1: 1,3",
    "2: 2,0",
    "3: 0,3",
    "4: 1,3",
    "5: 2,3"
)

